$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlWhole = 1

function Get-RowByLabel($label) {
    # LookAt:=xlWhole ensures an exact match (e.g. "RM 2" must not match "RM 232")
    $found = $ws.Range("A1:A1000").Find($label, [Type]::Missing, [Type]::Missing, $xlWhole)
    if ($found -eq $null) {
        throw "Label not found: $label"
    }
    return $found.Row
}

# --- Remove the two rows that were dropped from the source data ---
$rowToDelete = Get-RowByLabel "RM 232"
$ws.Rows($rowToDelete).Delete()

$rowToDelete = Get-RowByLabel "SC 92"
$ws.Rows($rowToDelete).Delete()

# --- Apply the per-row cell value corrections (rows identified by their ID in column A) ---

$r = Get-RowByLabel "RM 2"
$ws.Cells.Item($r, 3).Value = 14.9

$r = Get-RowByLabel "RM 8"
$ws.Cells.Item($r, 6).ClearContents()

$r = Get-RowByLabel "RM 9"
$ws.Cells.Item($r, 6).Value = 17.97

$r = Get-RowByLabel "RM 14"
$ws.Cells.Item($r, 6).ClearContents()

$r = Get-RowByLabel "RM 21"
$ws.Cells.Item($r, 3).ClearContents()

$r = Get-RowByLabel "RM 38"
$ws.Cells.Item($r, 6).ClearContents()

$r = Get-RowByLabel "RM 81"
$ws.Cells.Item($r, 3).Value = 12.5

$r = Get-RowByLabel "RM 90"
$ws.Cells.Item($r, 3).ClearContents()

$r = Get-RowByLabel "RM 95"
$ws.Cells.Item($r, 6).Value = 16.2

$r = Get-RowByLabel "RM 120"
$ws.Cells.Item($r, 6).Value = 18.35

$r = Get-RowByLabel "RM 125"
$ws.Cells.Item($r, 6).ClearContents()

$r = Get-RowByLabel "RM 134"
$ws.Cells.Item($r, 3).Value = 12.5

$r = Get-RowByLabel "RM 135"
$ws.Cells.Item($r, 3).Value = 12.7

$r = Get-RowByLabel "RM 138"
$ws.Cells.Item($r, 6).ClearContents()

$r = Get-RowByLabel "RM 140"
$ws.Cells.Item($r, 3).ClearContents()
$ws.Cells.Item($r, 6).Value = 16.48

$r = Get-RowByLabel "RM 142a"
$ws.Cells.Item($r, 3).ClearContents()

$r = Get-RowByLabel "RM 145"
$ws.Cells.Item($r, 6).Value = 16.6

$r = Get-RowByLabel "SC 5"
$ws.Cells.Item($r, 2).Value = -20.2

$r = Get-RowByLabel "SC 101"
$ws.Cells.Item($r, 2).ClearContents()
$ws.Cells.Item($r, 6).ClearContents()

$r = Get-RowByLabel "SC 120"
$ws.Cells.Item($r, 2).Value = -19.7

$r = Get-RowByLabel "SC 132"
$ws.Cells.Item($r, 3).Value = 15.3

$r = Get-RowByLabel "SC 193"
$ws.Cells.Item($r, 2).ClearContents()

$r = Get-RowByLabel "SC 232"
$ws.Cells.Item($r, 3).Value = 10.4
